$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values per row (columns I and K are dropped from the results in this run)
$rowData = @{
    2 = @{ 'C'=0.993979075244931; 'D'=0.9949798479784844; 'E'=0.9973540890848477; 'F'=1; 'G'=0.9829071632654329; 'H'=0.9958833301782321; 'J'=1; 'L'=0.9936388026113527; 'M'=0.9937976562499454; 'N'=0.9970152127389644; 'O'=0.9793297702495547 }
    3 = @{ 'C'=0.993624471265469; 'D'=0.9950699066052769; 'E'=0.9974572912098776; 'F'=1; 'G'=0.9839461463727678; 'H'=0.9959734705819233; 'J'=1; 'L'=0.9933313763553113; 'M'=0.9939720253226711; 'N'=0.9971654777130754; 'O'=0.9806699053221377 }
    4 = @{ 'C'=0.9949192836467563; 'D'=0.9945438180640228; 'E'=0.9971231989139545; 'F'=1; 'G'=0.9788562563145931; 'H'=0.9954469043308607; 'J'=1; 'L'=0.9946266233357345; 'M'=0.9930417482285482; 'N'=0.9968312743272536; 'O'=0.9740983416572403 }
    5 = @{ 'C'=0.9948378630038087; 'D'=0.994592133159213; 'E'=0.9971538814777505; 'F'=0.9999999999999999; 'G'=0.9792662085892417; 'H'=0.9954952632981243; 'J'=0.9999999999999998; 'L'=0.9945451754158078; 'M'=0.9931219504702389; 'N'=0.9968619670981913; 'O'=0.9746285535069026 }
    6 = @{ 'C'=0.9948144376853374; 'D'=0.9944356278807013; 'E'=0.9961851265146134; 'F'=1; 'G'=0.9790718602697165; 'H'=0.9953386159064495; 'J'=1; 'L'=0.9934539069838949; 'M'=0.9929621250551242; 'N'=0.9948280440691954; 'O'=0.974421573129465 }
    7 = @{ 'C'=0.9917838238992799; 'D'=0.9953798144180316; 'E'=0.9959893392176188; 'F'=1; 'G'=0.9899289126363155; 'H'=0.9962836598035857; 'J'=0.9999999999999999; 'L'=0.989212288900464; 'M'=0.9947864902534364; 'N'=0.99344047531562; 'O'=0.9883517816410423 }
    8 = @{ 'C'=0.992189047554438; 'D'=0.9952458554843827; 'E'=0.9956501536414539; 'F'=0.9999999999999999; 'G'=0.9887075956889463; 'H'=0.9961495792297728; 'J'=1; 'L'=0.9892075475422283; 'M'=0.9945560078571282; 'N'=0.9926908638383636; 'O'=0.9868209137500574 }
    9 = @{ 'C'=0.992194763370253; 'D'=0.9953111323942189; 'E'=0.9959464785928553; 'F'=0.9999999999999999; 'G'=0.9888875773188175; 'H'=0.9962149154137; 'J'=0.9999999999999997; 'L'=0.9896254591078131; 'M'=0.9946314714659704; 'N'=0.9933973853939897; 'O'=0.9870339780230196 }
    10 = @{ 'C'=0.9942931649584107; 'D'=0.9945935665439163; 'E'=0.9958739819436905; 'F'=0.9999999999999999; 'G'=0.9811286363385037; 'H'=0.995496697984399; 'J'=0.9999999999999998; 'L'=0.9923310422161266; 'M'=0.9932879257742934; 'N'=0.9939180536628169; 'O'=0.9771023998958911 }
    11 = @{ 'C'=0.9945206118438148; 'D'=0.9943186432622969; 'E'=0.9956526852566987; 'F'=1; 'G'=0.979198060527379; 'H'=0.9952215250612498; 'J'=1; 'L'=0.9924889730402187; 'M'=0.9928666145481658; 'N'=0.9936256733170646; 'O'=0.9746283578002397 }
    12 = @{ 'C'=0.993083787769835; 'D'=0.9950733431778213; 'E'=0.9958207131780678; 'F'=0.9999999999999998; 'G'=0.9860941143927945; 'H'=0.9959769102750152; 'J'=1; 'L'=0.9905588212844966; 'M'=0.9941662994377274; 'N'=0.993310203424481; 'O'=0.9834825857359527 }
    13 = @{ 'C'=0.9912600423546456; 'D'=0.9953689591844146; 'E'=0.9958188919924562; 'F'=1; 'G'=0.9908171618147896; 'H'=0.996272794712975; 'J'=0.9999999999999998; 'L'=0.9884162843909573; 'M'=0.9948552458847431; 'N'=0.99300278905879; 'O'=0.9894885316301043 }
    14 = @{ 'C'=0.9910473724789013; 'D'=0.9951476542302997; 'E'=0.9948155183826323; 'F'=1; 'G'=0.9905052367312093; 'H'=0.9960512888049539; 'J'=0.9999999999999999; 'L'=0.9867338485001819; 'M'=0.9946259531284898; 'N'=0.9905391020553729; 'O'=0.9891519136333979 }
    15 = @{ 'C'=0.9910122741661758; 'D'=0.9953816720402139; 'E'=0.9958997106843204; 'F'=1; 'G'=0.991228815928906; 'H'=0.9962855191125647; 'J'=1; 'L'=0.9882829098628939; 'M'=0.9949031959108006; 'N'=0.993198593533112; 'O'=0.9900080232009912 }
    16 = @{ 'C'=0.9898648861266423; 'D'=0.995418453217066; 'E'=0.9964109003453668; 'F'=1; 'G'=0.9928663990099267; 'H'=0.9963223336882226; 'J'=0.9999999999999999; 'L'=0.9878993229734896; 'M'=0.9950806251766237; 'N'=0.9944709492595286; 'O'=0.9920607474815649 }
    17 = @{ 'C'=0.9892884398954911; 'D'=0.9953851572537726; 'E'=0.996485320239295; 'F'=0.9999999999999998; 'G'=0.9934998378097286; 'H'=0.9962890074908392; 'J'=0.9999999999999999; 'L'=0.9874660410856098; 'M'=0.9951054367543352; 'N'=0.994688630587452; 'O'=0.9928573187592444 }
    18 = @{ 'C'=0.9940477591681391; 'D'=0.9947691887379668; 'E'=0.9960381552244415; 'F'=1; 'G'=0.9825029202106044; 'H'=0.9956724796505496; 'J'=0.9999999999999998; 'L'=0.9921608821781657; 'M'=0.9935698323709092; 'N'=0.9941587177363297; 'O'=0.9788635525865069 }
    19 = @{ 'C'=0.9929467789631036; 'D'=0.9906928642665799; 'E'=0.9921593072185186; 'F'=1; 'G'=0.9611089559593372; 'H'=0.9915924537106279; 'J'=0.9999999999999999; 'L'=0.9889815787837511; 'M'=0.987994548666092; 'N'=0.9881870229875505; 'O'=0.9514546711419518 }
    20 = @{ 'C'=0.9857287020906186; 'D'=0.9880009352247359; 'E'=0.9866423475669851; 'F'=1; 'G'=0.9593582825552284; 'H'=0.9888980802876466; 'J'=1; 'L'=0.9744540606998332; 'M'=0.9853821903594211; 'N'=0.9753947335359675; 'O'=0.9500463009699752 }
    21 = @{ 'C'=0.9887036047845105; 'D'=0.9889753876248141; 'E'=0.9888106941146061; 'F'=0.9999999999999999; 'G'=0.9593023735143466; 'H'=0.9898734175301662; 'J'=1; 'L'=0.9805567510029491; 'M'=0.9862699954633181; 'N'=0.9806660341108467; 'O'=0.9496131417852303 }
    22 = @{ 'C'=0.9922265071497063; 'D'=0.9896294137836643; 'E'=0.9913736146325236; 'F'=1; 'G'=0.9559664842551042; 'H'=0.9905280375713889; 'J'=1; 'L'=0.9880763395477189; 'M'=0.9865929695426309; 'N'=0.9872153530177266; 'O'=0.9447938646892488 }
    23 = @{ 'C'=0.9941299732920263; 'D'=0.9904691482753921; 'E'=0.992959334455772; 'F'=1; 'G'=0.957153793920178; 'H'=0.9913685345762143; 'J'=1; 'L'=0.9917270766103421; 'M'=0.9874610559346365; 'N'=0.9905505710266135; 'O'=0.9461128576079337 }
    24 = @{ 'C'=0.995247532258513; 'D'=0.9906197926260227; 'E'=0.9938711278553529; 'F'=1; 'G'=0.9556219444886158; 'H'=0.991519315718047; 'J'=1; 'L'=0.9940435336590371; 'M'=0.9874737522007472; 'N'=0.9926641372389237; 'O'=0.9439478778856533 }
    25 = @{ 'C'=0.9958736075553308; 'D'=0.9913165407382651; 'E'=0.9946920874470705; 'F'=0.9999999999999999; 'G'=0.9585123290777529; 'H'=0.9922166965059341; 'J'=0.9999999999999999; 'L'=0.9951561264273685; 'M'=0.9883503458902334; 'N'=0.9939732705840985; 'O'=0.9476640323827289 }
}

foreach ($row in $rowData.Keys) {
    $cols = $rowData[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
    # Columns I and K no longer have data for this timestep
    $ws.Range("I$row").ClearContents()
    $ws.Range("K$row").ClearContents()
}
